# Update cryptocurrency price/volume data in the sheet.
# Each entry is forced to Text format before assignment so that
# numeric-looking strings (e.g. "5.030", "18.10", "29.261.83") are
# preserved exactly as text rather than being coerced into doubles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "29.261.83" },
    @{ Cell = "E2"; Value = "  +0.40%  " },
    @{ Cell = "D3"; Value = "1.844.84" },
    @{ Cell = "E3"; Value = "  +0.52%  " },
    @{ Cell = "E4"; Value = "  +0.03%  " },
    @{ Cell = "D5"; Value = "242.41" },
    @{ Cell = "E5"; Value = "  +0.68%  " },
    @{ Cell = "D6"; Value = "0.6645" },
    @{ Cell = "E6"; Value = "  -0.10%  " },
    @{ Cell = "D7"; Value = "0.9997" },
    @{ Cell = "E7"; Value = "  +0.01%  " },
    @{ Cell = "D8"; Value = "0.07467" },
    @{ Cell = "E8"; Value = "  +1.46%  " },
    @{ Cell = "D9"; Value = "0.2956" },
    @{ Cell = "E9"; Value = "  +0.34%  " },
    @{ Cell = "D10"; Value = "23.46" },
    @{ Cell = "E10"; Value = "  +2.99%  " },
    @{ Cell = "D11"; Value = "0.07768" },
    @{ Cell = "E11"; Value = "  +0.91%  " },
    @{ Cell = "D12"; Value = "1.845.54" },
    @{ Cell = "E12"; Value = "  -0.60%  " },
    @{ Cell = "D13"; Value = "5.030" },
    @{ Cell = "E13"; Value = "  +0.29%  " },
    @{ Cell = "D14"; Value = "0.6744" },
    @{ Cell = "E14"; Value = "  +0.11%  " },
    @{ Cell = "D15"; Value = "83.66" },
    @{ Cell = "E15"; Value = "  -2.90%  " },
    @{ Cell = "D16"; Value = "6.203" },
    @{ Cell = "E16"; Value = "  +0.15%  " },
    @{ Cell = "D17"; Value = "0.000008631" },
    @{ Cell = "E17"; Value = "  +5.06%  " },
    @{ Cell = "D18"; Value = "29.265.56" },
    @{ Cell = "E18"; Value = "  +1.18%  " },
    @{ Cell = "D19"; Value = "2.096.97" },
    @{ Cell = "E19"; Value = "  +1.55%  " },
    @{ Cell = "D20"; Value = "228.66" },
    @{ Cell = "E20"; Value = "  +0.27%  " },
    @{ Cell = "D21"; Value = "12.56" },
    @{ Cell = "E21"; Value = "  +0.51%  " },
    @{ Cell = "D22"; Value = "0.9999" },
    @{ Cell = "E22"; Value = "  +0.07%  " },
    @{ Cell = "D23"; Value = "7.218" },
    @{ Cell = "E23"; Value = "  -0.42%  " },
    @{ Cell = "D24"; Value = "0.9994" },
    @{ Cell = "D25"; Value = "159.09" },
    @{ Cell = "E25"; Value = "  -0.71%  " },
    @{ Cell = "D26"; Value = "0.1412" },
    @{ Cell = "E26"; Value = "  +1.11%  " },
    @{ Cell = "D27"; Value = "8.652" },
    @{ Cell = "E27"; Value = "  -0.52%  " },
    @{ Cell = "D28"; Value = "18.10" },
    @{ Cell = "E28"; Value = "  +0.66%  " },
    @{ Cell = "D29"; Value = "1.510" },
    @{ Cell = "E29"; Value = "  +0.42%  " },
    @{ Cell = "D30"; Value = "4.137" },
    @{ Cell = "E30"; Value = "  -1.49%  " },
    @{ Cell = "D31"; Value = "4.063" },
    @{ Cell = "E31"; Value = "  -0.21%  " },
    @{ Cell = "D32"; Value = "1.193" },
    @{ Cell = "E32"; Value = "  +0.50%  " },
    @{ Cell = "D33"; Value = "0.05346" },
    @{ Cell = "E33"; Value = "  +0.41%  " },
    @{ Cell = "D34"; Value = "1.889" },
    @{ Cell = "E34"; Value = "  +1.55%  " },
    @{ Cell = "D35"; Value = "0.7459" },
    @{ Cell = "E35"; Value = "  -0.60%  " },
    @{ Cell = "D36"; Value = "1.157" },
    @{ Cell = "E36"; Value = "  +2.22%  " },
    @{ Cell = "D37"; Value = "2.650" },
    @{ Cell = "E37"; Value = "  -1.05%  " },
    @{ Cell = "D38"; Value = "1.325.06" },
    @{ Cell = "E38"; Value = "  +0.34%  " },
    @{ Cell = "D39"; Value = "0.01801" },
    @{ Cell = "E39"; Value = "  -0.21%  " },
    @{ Cell = "D40"; Value = "2.751" },
    @{ Cell = "E40"; Value = "  +0.75%  " },
    @{ Cell = "D41"; Value = "6.414" },
    @{ Cell = "E41"; Value = "  +7.43%  " },
    @{ Cell = "D42"; Value = "0.9187" },
    @{ Cell = "E42"; Value = "  -0.27%  " },
    @{ Cell = "D43"; Value = "0.9992" },
    @{ Cell = "E43"; Value = "  +0.05%  " },
    @{ Cell = "D44"; Value = "103.30" },
    @{ Cell = "E44"; Value = "  -0.19%  " },
    @{ Cell = "D45"; Value = "66.31" },
    @{ Cell = "E45"; Value = "  +3.67%  " },
    @{ Cell = "D46"; Value = "1.999.62" },
    @{ Cell = "E46"; Value = "  +1.85%  " },
    @{ Cell = "D47"; Value = "0.00000000123" },
    @{ Cell = "E47"; Value = "  -1.41%  " },
    @{ Cell = "D48"; Value = "0.5135" },
    @{ Cell = "D49"; Value = "0.07694" },
    @{ Cell = "E49"; Value = "  -5.22%  " },
    @{ Cell = "E50"; Value = "  +0.26%  " },
    @{ Cell = "D51"; Value = "0.05859" },
    @{ Cell = "E51"; Value = "  -1.19%  " }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}
